$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-shuffle existing region/color/reading values within same-date groups (rows 27-113) ---
$ws.Range("B27").Value = 1305; $ws.Range("C27").Value = "SOUTH EAST ASIA"
$ws.Range("B28").Value = 1374; $ws.Range("C28").Value = "WEST AFRICA"
$ws.Range("B29").Value = 1406; $ws.Range("C29").Value = "NORTH SEA"
$ws.Range("B30").Value = 1408; $ws.Range("C30").Value = "MENAM"; $ws.Range("D30").Value = "YELLOW"
$ws.Range("B31").Value = 1454; $ws.Range("C31").Value = "INDIA"; $ws.Range("D31").Value = "RED"
$ws.Range("B33").Value = 1602; $ws.Range("C33").Value = "MENAM"
$ws.Range("B34").Value = 1603; $ws.Range("C34").Value = "SOUTH EAST ASIA"
$ws.Range("B35").Value = 1709; $ws.Range("C35").Value = "WEST AFRICA"
$ws.Range("B36").Value = 1711; $ws.Range("C36").Value = "NORTH SEA"
$ws.Range("B37").Value = 1889; $ws.Range("C37").Value = "SOUTH EAST ASIA"
$ws.Range("B38").Value = 1891; $ws.Range("C38").Value = "INDIA"
$ws.Range("B39").Value = 1898; $ws.Range("C39").Value = "MENAM"; $ws.Range("D39").Value = "YELLOW"
$ws.Range("B40").Value = 1990; $ws.Range("C40").Value = "NORTH SEA"
$ws.Range("B41").Value = 1996; $ws.Range("C41").Value = "WEST AFRICA"; $ws.Range("D41").Value = "WHITE"
$ws.Range("B42").Value = 2138; $ws.Range("C42").Value = "MENAM"
$ws.Range("B43").Value = 2139; $ws.Range("C43").Value = "INDIA"
$ws.Range("B44").Value = 2152; $ws.Range("C44").Value = "SOUTH EAST ASIA"
$ws.Range("B45").Value = 2201; $ws.Range("C45").Value = "WEST AFRICA"
$ws.Range("B46").Value = 2249; $ws.Range("C46").Value = "NORTH SEA"
$ws.Range("B47").Value = 2417; $ws.Range("C47").Value = "SOUTH EAST ASIA"
$ws.Range("B51").Value = 2389; $ws.Range("C51").Value = "INDIA"
$ws.Range("B52").Value = 2765; $ws.Range("C52").Value = "MENAM"; $ws.Range("D52").Value = "YELLOW"
$ws.Range("B53").Value = 2660; $ws.Range("C53").Value = "INDIA"
$ws.Range("B54").Value = 2748; $ws.Range("C54").Value = "NORTH SEA"; $ws.Range("D54").Value = "WHITE"
$ws.Range("B56").Value = 2703; $ws.Range("C56").Value = "WEST AFRICA"
$ws.Range("B57").Value = 2992; $ws.Range("C57").Value = "NORTH SEA"
$ws.Range("B60").Value = 2903; $ws.Range("C60").Value = "MENAM"
$ws.Range("B63").Value = 3298; $ws.Range("C63").Value = "NORTH SEA"; $ws.Range("D63").Value = "WHITE"
$ws.Range("B65").Value = 3336; $ws.Range("C65").Value = "SOUTH EAST ASIA"
$ws.Range("B66").Value = 3337; $ws.Range("C66").Value = "INDIA"; $ws.Range("D66").Value = "YELLOW"
$ws.Range("B67").Value = 3547; $ws.Range("C67").Value = "NORTH SEA"; $ws.Range("D67").Value = "WHITE"
$ws.Range("B68").Value = 3636; $ws.Range("C68").Value = "SOUTH EAST ASIA"; $ws.Range("D68").Value = "WHITE"
$ws.Range("B69").Value = 3462; $ws.Range("C69").Value = "WEST AFRICA"
$ws.Range("B70").Value = 3517; $ws.Range("C70").Value = "INDIA"; $ws.Range("D70").Value = "YELLOW"
$ws.Range("B71").Value = 3411; $ws.Range("C71").Value = "MENAM"; $ws.Range("D71").Value = "YELLOW"
$ws.Range("B72").Value = 3784; $ws.Range("C72").Value = "MENAM"; $ws.Range("D72").Value = "YELLOW"
$ws.Range("B73").Value = 3881; $ws.Range("C73").Value = "SOUTH EAST ASIA"; $ws.Range("D73").Value = "WHITE"
$ws.Range("B74").Value = 3827; $ws.Range("C74").Value = "NORTH SEA"
$ws.Range("B75").Value = 3737; $ws.Range("C75").Value = "WEST AFRICA"
$ws.Range("B77").Value = 4094; $ws.Range("C77").Value = "NORTH SEA"; $ws.Range("D77").Value = "WHITE"
$ws.Range("B79").Value = 4011; $ws.Range("C79").Value = "WEST AFRICA"
$ws.Range("B80").Value = 3966; $ws.Range("C80").Value = "MENAM"; $ws.Range("D80").Value = "YELLOW"
$ws.Range("B81").Value = 4186; $ws.Range("C81").Value = "SOUTH EAST ASIA"
$ws.Range("B82").Value = 4217; $ws.Range("C82").Value = "MENAM"
$ws.Range("B83").Value = 4311; $ws.Range("C83").Value = "INDIA"; $ws.Range("D83").Value = "YELLOW"
$ws.Range("B84").Value = 4272; $ws.Range("C84").Value = "WEST AFRICA"; $ws.Range("D84").Value = "WHITE"
$ws.Range("B85").Value = 4397; $ws.Range("C85").Value = "NORTH SEA"
$ws.Range("B86").Value = 4540; $ws.Range("C86").Value = "WEST AFRICA"
$ws.Range("B87").Value = 4679; $ws.Range("C87").Value = "NORTH SEA"
$ws.Range("B90").Value = 4910; $ws.Range("C90").Value = "NORTH SEA"
$ws.Range("B91").Value = 4872; $ws.Range("C91").Value = "SOUTH EAST ASIA"
$ws.Range("B92").Value = 4828; $ws.Range("C92").Value = "WEST AFRICA"; $ws.Range("D92").Value = "WHITE"
$ws.Range("B93").Value = 4877; $ws.Range("C93").Value = "INDIA"
$ws.Range("B94").Value = 4807; $ws.Range("C94").Value = "MENAM"; $ws.Range("D94").Value = "YELLOW"
$ws.Range("B105").Value = 5766; $ws.Range("C105").Value = "NORTH SEA"; $ws.Range("D105").Value = "WHITE"
$ws.Range("B108").Value = 5742; $ws.Range("C108").Value = "INDIA"; $ws.Range("D108").Value = "YELLOW"
$ws.Range("B109").Value = 6042; $ws.Range("C109").Value = "MENAM"
$ws.Range("B110").Value = 6003; $ws.Range("C110").Value = "INDIA"
$ws.Range("B111").Value = 5998; $ws.Range("C111").Value = "WEST AFRICA"; $ws.Range("D111").Value = "WHITE"
$ws.Range("B112").Value = 6077; $ws.Range("C112").Value = "NORTH SEA"
$ws.Range("B113").Value = 6035; $ws.Range("C113").Value = "SOUTH EAST ASIA"; $ws.Range("D113").Value = "YELLOW"

# --- Append new rows 114-119 (date 2023-04-05 / serial 45021) ---
# Copy formatting (style indexes) from the last existing data row (113) first
$ws.Range("A113").Copy() | Out-Null
$ws.Range("A114:A119").PasteSpecial(-4122) | Out-Null
$ws.Range("E113").Copy() | Out-Null
$ws.Range("E114:E119").PasteSpecial(-4122) | Out-Null

# row 114
$ws.Range("A114").Value = 112
$ws.Range("B114").Value = 6348
$ws.Range("C114").Value = "MENAM"
$ws.Range("D114").Value = "YELLOW"
$ws.Range("E114").Value = 45021

# row 115
$ws.Range("A115").Value = 113
$ws.Range("B115").Value = 6302
$ws.Range("C115").Value = "SOUTH EAST ASIA"
$ws.Range("D115").Value = "WHITE"
$ws.Range("E115").Value = 45021

# row 116
$ws.Range("A116").Value = 114
$ws.Range("B116").Value = 6253
$ws.Range("C116").Value = "NORTH SEA"
$ws.Range("D116").Value = "WHITE"
$ws.Range("E116").Value = 45021

# row 117
$ws.Range("A117").Value = 115
$ws.Range("B117").Value = 6254
$ws.Range("C117").Value = "WEST AFRICA"
$ws.Range("D117").Value = "WHITE"
$ws.Range("E117").Value = 45021

# row 118
$ws.Range("A118").Value = 116
$ws.Range("B118").Value = 6276
$ws.Range("C118").Value = "INDIA"
$ws.Range("D118").Value = "WHITE"
$ws.Range("E118").Value = 45021

# row 119
$ws.Range("A119").Value = 117
$ws.Range("B119").Value = 6546
$ws.Range("D119").Value = "WHITE"
$ws.Range("E119").Value = 45021

